# Add a new "Thumbnail URL" column to the books sheet, right before the
# existing "Shelf-A" column. Inserting a whole column shifts the former
# H:L ("Shelf-A".."Shelf-E") one slot to the right (I:M) and Excel grows
# the sheet's dimension/used-range accordingly (A1:L1 -> A1:M1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at H (shifts Shelf-A..Shelf-E from H:L to I:M).
$ws.Columns.Item(8).EntireColumn.Insert()

# Populate the new header cell.
$ws.Cells.Item(1, 8).Value = "Thumbnail URL"
